$d = $word.ActiveDocument

# 1. "Great news! We will be in [CITY NAME] from" -> "¡Buenas noticias! We will be in [CITY NAME] from"
$r = $d.Content
$r.Find.Execute("Great news! We will be in [CITY NAME] from", $true, $false, $false, $false, $false, $true, 1, $false, "¡Buenas noticias! We will be in [CITY NAME] from", 2)

# 2. "P.S. We’re giving out free Deriv merchandise. Don’t miss out!" -> "P.S. We’re giving out free Deriv merchandise. ¡No se lo pierda!"
$r = $d.Content
$r.Find.Execute("P.S. We’re giving out free Deriv merchandise. Don’t miss out!", $true, $false, $false, $false, $false, $true, 1, $false, "P.S. We’re giving out free Deriv merchandise. ¡No se lo pierda!", 2)

# 3. "If you have questions, contact us " -> "Si tiene alguna pregunta, póngase en contacto con nosotros "
$r = $d.Content
$r.Find.Execute("If you have questions, contact us ", $true, $false, $false, $false, $false, $true, 1, $false, "Si tiene alguna pregunta, póngase en contacto con nosotros ", 2)
